$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 106. This shifts the existing rows
# 106-131 down to 107-132 (formatting included), matching the target
# workbook where the old row106..row131 data now lives in row107..row132
# and a brand new row106 is introduced.
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the new record.
$ws.Cells.Item(106,1).Value  = 8
$ws.Cells.Item(106,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(106,3).Value  = "Coquimbo"
$ws.Cells.Item(106,4).Value  = [DateTime]"2023-01-06"
$ws.Cells.Item(106,5).Value  = 4
$ws.Cells.Item(106,6).Value  = "Fruta"
$ws.Cells.Item(106,7).Value  = 100109
$ws.Cells.Item(106,8).Value  = "Uva"
$ws.Cells.Item(106,9).Value  = 100109001
$ws.Cells.Item(106,10).Value = "Uva"
$ws.Cells.Item(106,11).Value = "Flame Seedless"
$ws.Cells.Item(106,12).Value = "Primera"
$ws.Cells.Item(106,13).Value = 1100
$ws.Cells.Item(106,14).Value = 9000
$ws.Cells.Item(106,15).Value = 10000
$ws.Cells.Item(106,16).Value = 9500
$ws.Cells.Item(106,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(106,18).Value = "Provincia de Limarí"
$ws.Cells.Item(106,19).Value = 950
$ws.Cells.Item(106,20).Value = 10

# Ensure the date cell keeps the expected date/time number format used by
# the rest of the "D" column.
$ws.Cells.Item(106,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
